$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 71431496
$ws.Range("J40").Value = 250001460
$ws.Range("L40").Value = 250001460
$ws.Range("N40").Value = -250001810

$ws.Range("H76").Value = 7874.2
$ws.Range("J76").Value = 4394.5
$ws.Range("L76").Value = 4394.5
$ws.Range("N76").Value = -5024.5

$ws.Range("H79").Value = 7874.2
$ws.Range("J79").Value = 4394.5
$ws.Range("L79").Value = 4394.5
$ws.Range("N79").Value = -6578.5

$ws.Range("H96").Value = 1451231.8
$ws.Range("I96").Value = 1723.7142
$ws.Range("J96").Value = 4833417.5
$ws.Range("K96").Value = 5171.142599999999
$ws.Range("L96").Value = 14500252.5
$ws.Range("M96").Value = -3798.142599999999
$ws.Range("N96").Value = -14502998.5

$ws.Range("H98").Value = 10418988
$ws.Range("I98").Value = 11365714
$ws.Range("K98").Value = 11365714
$ws.Range("M98").Value = -11364216

$ws.Range("H100").Value = 10809.637
$ws.Range("I100").Value = 5961.6
$ws.Range("K100").Value = 5961.6
$ws.Range("M100").Value = -5420.6

$ws.Range("H122").Value = 10418988
$ws.Range("I122").Value = 11365714
$ws.Range("K122").Value = 34097142
$ws.Range("M122").Value = -34094692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3896.2954
$ws.Range("I32").Value = 3545.0466
$ws.Range("K32").Value = 3545.0466
$ws.Range("M32").Value = -3258.0466

$ws.Range("H97").Value = 1274.2307
$ws.Range("I97").Value = 1108.0952
$ws.Range("J97").Value = 1972
$ws.Range("K97").Value = 1108.0952
$ws.Range("L97").Value = 1972
$ws.Range("M97").Value = -612.0952
$ws.Range("N97").Value = -2964

$ws.Range("H102").Value = 1775.4445
$ws.Range("I102").Value = 1161.3334
$ws.Range("J102").Value = 3003.6667
$ws.Range("K102").Value = 1161.3334
$ws.Range("L102").Value = 3003.6667
$ws.Range("M102").Value = 460.6666
$ws.Range("N102").Value = -6247.6667

$ws.Range("H132").Value = 3034683.8
$ws.Range("I132").Value = 4259.769
$ws.Range("J132").Value = 14290545
$ws.Range("K132").Value = 12779.307
$ws.Range("L132").Value = 42871635
$ws.Range("M132").Value = -10249.307
$ws.Range("N132").Value = -42876695

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1719.8125
$ws.Range("I94").Value = 1766.3043
$ws.Range("K94").Value = 1766.3043
$ws.Range("M94").Value = -1315.3043

$ws.Range("H134").Value = 16670844
$ws.Range("I134").Value = 3762.25
$ws.Range("K134").Value = 11286.75
$ws.Range("M134").Value = -8751.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H86").Value = 11155.714
$ws.Range("I86").Value = 14165.889
$ws.Range("K86").Value = 14165.889
$ws.Range("M86").Value = -13042.889

$ws.Range("H89").Value = 11155.714
$ws.Range("I89").Value = 14165.889
$ws.Range("K89").Value = 70829.44499999999
$ws.Range("M89").Value = -65213.44499999999

$ws.Range("H99").Value = 9599.6
$ws.Range("J99").Value = 6500
$ws.Range("L99").Value = 6500
$ws.Range("N99").Value = -9496

$ws.Range("H126").Value = 9599.6
$ws.Range("J126").Value = 6500
$ws.Range("L126").Value = 19500
$ws.Range("N126").Value = -24440

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.833332
$ws.Range("I2").Value = 40.5
$ws.Range("J2").Value = 56
$ws.Range("K2").Value = 243
$ws.Range("L2").Value = 336
$ws.Range("M2").Value = -130
$ws.Range("N2").Value = -562

$ws.Range("H4").Value = 4092635.5
$ws.Range("I4").Value = 3248726.5
$ws.Range("J4").Value = 10000000
$ws.Range("K4").Value = 9746179.5
$ws.Range("L4").Value = 30000000
$ws.Range("M4").Value = -9746067.5
$ws.Range("N4").Value = -30000224

$ws.Range("H39").Value = 12442.667
$ws.Range("J39").Value = 17664
$ws.Range("L39").Value = 52992
$ws.Range("N39").Value = -53580

$ws.Range("H50").Value = 963.3333
$ws.Range("I50").Value = 963.3333
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 2889.9999
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -2408.9999
$ws.Range("N50").ClearContents()

$ws.Range("H53").Value = 963.3333
$ws.Range("I53").Value = 963.3333
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 2889.9999
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -2408.9999
$ws.Range("N53").ClearContents()

$ws.Range("H121").Value = 4464.294
$ws.Range("I121").Value = 440
$ws.Range("J121").Value = 4715.8125
$ws.Range("K121").Value = 1320
$ws.Range("L121").Value = 14147.4375
$ws.Range("M121").Value = -10
$ws.Range("N121").Value = -16767.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3438.111
$ws.Range("J80").Value = 3566.4285
$ws.Range("L80").Value = 3566.4285
$ws.Range("N80").Value = -5562.4285

$ws.Range("H83").Value = 3438.111
$ws.Range("J83").Value = 3566.4285
$ws.Range("L83").Value = 17832.1425
$ws.Range("N83").Value = -27816.1425

$ws.Range("H97").Value = 853.3333
$ws.Range("I97").Value = 1082.7778
$ws.Range("J97").Value = 165
$ws.Range("K97").Value = 1082.7778
$ws.Range("L97").Value = 165
$ws.Range("M97").Value = -586.7778000000001
$ws.Range("N97").Value = -1157

$ws.Range("H122").Value = 2214.9
$ws.Range("I122").Value = 1792.4
$ws.Range("K122").Value = 5377.200000000001
$ws.Range("M122").Value = -2927.200000000001

$ws.Range("H132").Value = 5884611.5
$ws.Range("I132").Value = 1664.4445
$ws.Range("J132").Value = 12502927
$ws.Range("K132").Value = 4993.333500000001
$ws.Range("L132").Value = 37508781
$ws.Range("M132").Value = -2463.333500000001
$ws.Range("N132").Value = -37513841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 471.8
$ws.Range("I22").Value = 473.5
$ws.Range("K22").Value = 473.5
$ws.Range("M22").Value = -178.5

$ws.Range("H27").Value = 471.8
$ws.Range("I27").Value = 473.5
$ws.Range("K27").Value = 473.5
$ws.Range("M27").Value = -366.5

$ws.Range("H40").Value = 5058.8125
$ws.Range("I40").Value = 3919.3076
$ws.Range("K40").Value = 3919.3076
$ws.Range("M40").Value = -3783.3076

$ws.Range("H46").Value = 2700
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376

$ws.Range("H93").Value = 2226307.8
$ws.Range("J93").Value = 6180434.5
$ws.Range("L93").Value = 6180434.5
$ws.Range("N93").Value = -6182930.5

$ws.Range("H122").Value = 3422.122
$ws.Range("J122").Value = 5168.6665
$ws.Range("L122").Value = 15505.9995
$ws.Range("N122").Value = -20405.9995

$ws.Range("I132").Value = 2291.4736
$ws.Range("J132").Value = 4446.385
$ws.Range("K132").Value = 6874.4208
$ws.Range("L132").Value = 13339.155
$ws.Range("M132").Value = -4344.4208
$ws.Range("N132").Value = -18399.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 110999
$ws.Range("I2").Value = 110999
$ws.Range("K2").Value = 110999
$ws.Range("M2").Value = -110887

$ws.Range("H126").Value = 3566.4285
$ws.Range("I126").Value = 3860.8333
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 11582.4999
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -9112.499899999999
$ws.Range("N126").Value = -10340

$ws.Range("H132").Value = 460212.22
$ws.Range("I132").Value = 5538.8667
$ws.Range("J132").Value = 1434512.2
$ws.Range("K132").Value = 16616.6001
$ws.Range("L132").Value = 4303536.6
$ws.Range("M132").Value = -14086.6001
$ws.Range("N132").Value = -4308596.6
